$d = $word.ActiveDocument

$find = "για τον αστερισμό του Αστερισμός του Κύκνου"
$repl = "για τον  Αστερισμός του Κύκνου"

$rng = $d.Content
$rng.Find.ClearFormatting()
$rng.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $repl, 2)
